# Updated symbol list on Thu Dec 15 17:15:42 UTC 2022 with GitHub Actions
# This script applies the cell-value changes captured in the commit diff:
#   - "Hora" (column G) bumps from 16 -> 17 for every data row (2-51)
#   - Several "Price" (column D) values are refreshed with newer quotes
#   - Rows 10-18 shift: coin rankings 9-17 (WazirX, MandalaExchangeToken,
#     LiechtensteinCryptoassetsExchange, BitrueCoin, BitMartToken, MCDex,
#     BitForexToken, CoinExToken, One) move up one row with new data,
#     updating the Coin (B), Link (C) and Volume(1h) (E) columns as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose stored values look numeric ("Price" and "Hora" columns).
# Excel would otherwise silently reinterpret a numeric-looking string as a
# true number (dropping trailing zeros / changing precision), so each of
# these cells is explicitly formatted as Text before the value is written,
# keeping them identical in shape to the original inline strings.
$numericTextUpdates = @(
    @{ Cell = "D2"; Value = "263.22" },
    @{ Cell = "G2"; Value = "17" },
    @{ Cell = "D3"; Value = "22.79" },
    @{ Cell = "G3"; Value = "17" },
    @{ Cell = "D4"; Value = "6.201" },
    @{ Cell = "G4"; Value = "17" },
    @{ Cell = "D5"; Value = "0.06099" },
    @{ Cell = "G5"; Value = "17" },
    @{ Cell = "D6"; Value = "3.514" },
    @{ Cell = "G6"; Value = "17" },
    @{ Cell = "D7"; Value = "6.709" },
    @{ Cell = "G7"; Value = "17" },
    @{ Cell = "D8"; Value = "1.359" },
    @{ Cell = "G8"; Value = "17" },
    @{ Cell = "D9"; Value = "0.7980" },
    @{ Cell = "G9"; Value = "17" },
    @{ Cell = "D10"; Value = "0.1574" },
    @{ Cell = "G10"; Value = "17" },
    @{ Cell = "D11"; Value = "0.08136" },
    @{ Cell = "G11"; Value = "17" },
    @{ Cell = "D12"; Value = "0.03322" },
    @{ Cell = "G12"; Value = "17" },
    @{ Cell = "D13"; Value = "0.03156" },
    @{ Cell = "G13"; Value = "17" },
    @{ Cell = "D14"; Value = "0.09257" },
    @{ Cell = "G14"; Value = "17" },
    @{ Cell = "D15"; Value = "3.914" },
    @{ Cell = "G15"; Value = "17" },
    @{ Cell = "D16"; Value = "0.001694" },
    @{ Cell = "G16"; Value = "17" },
    @{ Cell = "D17"; Value = "0.04830" },
    @{ Cell = "G17"; Value = "17" },
    @{ Cell = "D18"; Value = "0.0006205" },
    @{ Cell = "G18"; Value = "17" },
    @{ Cell = "D19"; Value = "0.006193" },
    @{ Cell = "G19"; Value = "17" },
    @{ Cell = "D20"; Value = "0.001101" },
    @{ Cell = "G20"; Value = "17" },
    @{ Cell = "D21"; Value = "0.003190" },
    @{ Cell = "G21"; Value = "17" },
    @{ Cell = "G22"; Value = "17" },
    @{ Cell = "D23"; Value = "3.693" },
    @{ Cell = "G23"; Value = "17" },
    @{ Cell = "D24"; Value = "2.266" },
    @{ Cell = "G24"; Value = "17" },
    @{ Cell = "D25"; Value = "0.3384" },
    @{ Cell = "G25"; Value = "17" },
    @{ Cell = "G26"; Value = "17" },
    @{ Cell = "D27"; Value = "0.0004337" },
    @{ Cell = "G27"; Value = "17" },
    @{ Cell = "G28"; Value = "17" },
    @{ Cell = "G29"; Value = "17" },
    @{ Cell = "G30"; Value = "17" },
    @{ Cell = "G31"; Value = "17" },
    @{ Cell = "G32"; Value = "17" },
    @{ Cell = "G33"; Value = "17" },
    @{ Cell = "G34"; Value = "17" },
    @{ Cell = "G35"; Value = "17" },
    @{ Cell = "G36"; Value = "17" },
    @{ Cell = "G37"; Value = "17" },
    @{ Cell = "G38"; Value = "17" },
    @{ Cell = "G39"; Value = "17" },
    @{ Cell = "D40"; Value = "0.04596" },
    @{ Cell = "G40"; Value = "17" },
    @{ Cell = "D41"; Value = "0.007216" },
    @{ Cell = "G41"; Value = "17" },
    @{ Cell = "D42"; Value = "0.003904" },
    @{ Cell = "G42"; Value = "17" },
    @{ Cell = "G43"; Value = "17" },
    @{ Cell = "D44"; Value = "0.01093" },
    @{ Cell = "G44"; Value = "17" },
    @{ Cell = "D45"; Value = "0.002972" },
    @{ Cell = "G45"; Value = "17" },
    @{ Cell = "D46"; Value = "0.00006047" },
    @{ Cell = "G46"; Value = "17" },
    @{ Cell = "D47"; Value = "0.00000000751" },
    @{ Cell = "G47"; Value = "17" },
    @{ Cell = "D48"; Value = "0.7006" },
    @{ Cell = "G48"; Value = "17" },
    @{ Cell = "D49"; Value = "0.04850" },
    @{ Cell = "G49"; Value = "17" },
    @{ Cell = "D50"; Value = "0.00002102" },
    @{ Cell = "G50"; Value = "17" },
    @{ Cell = "D51"; Value = "0.01011" },
    @{ Cell = "G51"; Value = "17" }
)

foreach ($u in $numericTextUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

# Cells holding genuinely textual content (Coin name, Link, Volume(1h))
# can be assigned directly.
$textUpdates = @(
    @{ Cell = "B10"; Value = "WazirX" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "E10"; Value = "9WazirXWRX" },
    @{ Cell = "B11"; Value = "MandalaExchangeToken" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "E11"; Value = "10MandalaExchangeTokenMDX" },
    @{ Cell = "B12"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "E12"; Value = "11LiechtensteinCryptoassetsExchangeLCX" },
    @{ Cell = "B13"; Value = "BitrueCoin" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "E13"; Value = "12BitrueCoinBTR" },
    @{ Cell = "B14"; Value = "BitMartToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "E14"; Value = "13BitMartTokenBMX" },
    @{ Cell = "B15"; Value = "MCDex" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" },
    @{ Cell = "E15"; Value = "14MCDexMCB" },
    @{ Cell = "B16"; Value = "BitForexToken" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "E16"; Value = "15BitForexTokenBF" },
    @{ Cell = "B17"; Value = "CoinExToken" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" },
    @{ Cell = "E17"; Value = "16CoinExTokenCET" },
    @{ Cell = "B18"; Value = "One" },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" },
    @{ Cell = "E18"; Value = "17OneONEWorstin24h" },
    @{ Cell = "E49"; Value = "48BOLOBOLO" }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
